$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "10/01/2025"
$ws.Range("A45").ClearFormats()
$ws.Range("B45").Value = 14999.11
